# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 3.755628166162433)
    3  = @(0.6545652718822623, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 3.536033448013082)
    4  = @(0.2881169905109251, 109.9114832445916, 3.223369029078222, 13.86384647080068, 127.2868157349814)
    5  = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 13.86384647080068, 15.54439750930335)
    6  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    7  = @(1.445647641019636, 0.002658071450198252, 0.7210945179870265, 0.5333859586016987, 2.70278618905856)
    8  = @(1.445647641019636, 1.626987699542094, 189.6080260415259, 13.86384647080068, 206.5445078528883)
    9  = @(0.6545652718822623, 0.04103571897497393, 0.7210945179870265, 13.86384647080068, 15.28054197964495)
    10 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 16.29500630922404)
    11 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 24.14949828602258)
    12 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165)
    13 = @(0.1169995834814548, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1.412515779045154)
    14 = @(0.2881169905109251, 0.3048912486333797, 3.223369029078222, 13.86384647080068, 17.68022373902321)
    15 = @(0.04172184405617529, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.351702369198972)
    16 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    17 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 4.327115817150455)
    18 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    19 = @(0.2881169905109251, 1.626987699542094, 19575605.8673771, 14773364.14517103, 34348971.92765282)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
    $ws.Cells.Item($r, 7).Value = $vals[4]
}
